$wb = $excel.ActiveWorkbook

# Update "DateProd" (column B) timestamps on each sheet to reflect the
# latest Katalon test-run captures (BWP Object Repository refresh / VRelay
# fixes). Values applied in the chronological order they were produced.

$wsAuth = $wb.Worksheets.Item("CC-Payments-Auth")
$wsAuth.Range("B2").Value2 = "Thu Sep 04 06:23:30 IST 2025"
$wsAuth.Range("B3").Value2 = "Thu Sep 04 06:24:19 IST 2025"
$wsAuth.Range("B4").Value2 = "Thu Sep 04 06:25:05 IST 2025"
$wsAuth.Range("B5").Value2 = "Thu Sep 04 06:25:51 IST 2025"
$wsAuth.Range("B6").Value2 = "Thu Sep 04 06:26:34 IST 2025"
$wsAuth.Range("B7").Value2 = "Thu Sep 04 06:27:20 IST 2025"

$wsDebit = $wb.Worksheets.Item("ACH-Payments-Debit")
$wsDebit.Range("B2").Value2 = "Thu Sep 04 06:28:04 IST 2025"
$wsDebit.Range("B3").Value2 = "Thu Sep 04 06:28:48 IST 2025"
$wsDebit.Range("B4").Value2 = "Thu Sep 04 06:29:32 IST 2025"
$wsDebit.Range("B5").Value2 = "Thu Sep 04 06:30:15 IST 2025"
$wsDebit.Range("B6").Value2 = "Thu Sep 04 06:31:09 IST 2025"
$wsDebit.Range("B7").Value2 = "Thu Sep 04 06:31:54 IST 2025"
$wsDebit.Range("B8").Value2 = "Thu Sep 04 06:32:45 IST 2025"
$wsDebit.Range("B9").Value2 = "Thu Sep 04 06:33:38 IST 2025"
$wsDebit.Range("B10").Value2 = "Thu Sep 04 06:34:21 IST 2025"

$wsSale = $wb.Worksheets.Item("CC-Payments-Sale")
$wsSale.Range("B2").Value2 = "Thu Sep 04 06:35:07 IST 2025"
